# Help_ManExport.docx update:
# After the "Click the 'Export Flagged' button..." bullet, add a blank
# line and a new "Note: ... Marked cells will not be cleared ..." bullet.
# The existing `_GoBack` bookmark (which previously sat right after the
# "Export Flagged" run) is relocated into the new "Note:" paragraph,
# between the "Note: " run and the "Marked cells..." run.

$d = $word.ActiveDocument

# Find the paragraph that holds the "Export Flagged" bullet -- this is
# the last paragraph in the document body.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Export Flagged*") {
        $target = $p
    }
}

# Position right before that paragraph's end-of-paragraph mark.
$insertPos = $target.Range.End - 1

# The `_GoBack` bookmark currently sits at $insertPos (right after the
# "...will be exported. " text). Remove it here; it gets re-created
# further down, inside the newly inserted "Note:" paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$insertRange = $d.Range($insertPos, $insertPos)

$newParagraphsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Note: </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Marked cells will not be cleared after exporting. This will have to be done manually. </w:t></w:r></w:p>
'@

$insertRange.InsertXML($newParagraphsXml)
